# Atualiza notas dos alunos: preenche a coluna "TF" (H) com a nota da
# prova final e marca o "Conceito" (J) de cada aluno como aprovado ("A").
# A coluna "Nota" (I) é recalculada automaticamente pelas fórmulas SUM já
# existentes na planilha.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$notasFinais = @{
    2 = 2
    4 = 2.5
    5 = 2.5
    6 = 2.5
}

foreach ($linha in $notasFinais.Keys) {
    $ws.Cells.Item($linha, 8).Value = $notasFinais[$linha]   # coluna H (TF)
    $ws.Cells.Item($linha, 10).Value = "A"                   # coluna J (Conceito)
}

# A linha 3 já tinha a nota final preenchida; só falta o conceito.
$ws.Cells.Item(3, 10).Value = "A"

# Atualiza a célula selecionada, como ficou registrado no arquivo.
$ws.Range("H5").Select()
